# Updated ABS Ultra Durable profile
#
# Row 18 ("rigid.ink ABS Ultra Durable") previously answered "No, but
# would probably be helpful" for whether an enclosure was used, with no
# notes link. After testing the profile with an enclosure only, update
# the answer to "Yes" and add the "Notes" entry that the other "Yes"
# rows (17, 19, 21) already carry - copying their formatting so the new
# cell matches the rest of the column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bring F18's formatting in line with the other "Notes" cells in the
# column (F17/F19/F21) by copying one of them in, then overwrite the
# value with the correct text.
$ws.Range("F17").Copy()
$ws.Range("F18").PasteSpecial()
$ws.Range("F18").Value = "Notes"

# Enclosure answer changes from "No, but would probably be helpful" to
# "Yes" now that the profile has been validated with an enclosure.
$ws.Range("E18").Value = "Yes"

# The workbook was left with F18 selected when saved.
$ws.Range("F18").Select()
